# Update the "associate" regularNum sheet with the optimized program data.
# Only the numeric values in A2:J10 change; formatting/styles are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.01, 615, 426, 286, 194, 124, 61, 25, 9, 3)
    3  = @(0.02, 327, 219, 145, 95, 54, 23, 7, 2, 0)
    4  = @(0.03, 197, 142, 93, 59, 37, 13, 1, 0, 0)
    5  = @(0.04, 126, 99, 63, 40, 21, 10, 0, 0, 0)
    6  = @(0.05, 86, 69, 44, 29, 14, 6, 0, 0, 0)
    7  = @(0.06, 76, 62, 41, 27, 14, 6, 0, 0, 0)
    8  = @(0.07, 56, 49, 35, 23, 12, 5, 0, 0, 0)
    9  = @(0.08, 42, 38, 27, 21, 11, 4, 0, 0, 0)
    10 = @(0.09, 32, 31, 22, 17, 9, 3, 0, 0, 0)
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $values[$col - 1]
    }
}
